$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the set ordering / value swaps in column E ---
$ws.Range("E10").Value = "{'Tuple[NoneType]', 'Tuple[None]'}"
$ws.Range("E11").Value = "Tuple[NoneType]"
$ws.Range("E12").Value = "{'list', 'Tuple[NoneType]'}"
$ws.Range("E13").Value = "list"
$ws.Range("E16").Value = "{'callable', 'any'}"
$ws.Range("E17").Value = "callable"

# --- Restructure row 133: move label from C133 to E133, clear D133, set F133 as accuracy number ---
$ws.Range("C133").Value = ""
$ws.Range("D133").Value = ""
$ws.Range("E133").Value = "Scalpel Accuracy:"
$ws.Range("F133").Value = 93.84999999999999

# --- Update label text in row 134 ---
$ws.Range("E134").Value = "Accuracy vs PyType"
